$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Make sure the data range is formatted as Text so that numeric-looking
# values (e.g. "123") are stored as text, not numbers - matching the
# "Search element in table by 'Артикул'" feature which needs string
# comparisons across the whole table.
$range = $ws.Range("A2:F6")
$range.NumberFormat = "@"

$data = @(
    @("123",  "aaa",    "3",   "2",   "1",    "aaa"),
    @("1234", "fff",    "45",  "32",  "42",   "aaaasd"),
    @("231",  "asda",   "43",  "213", "23",   "aasda"),
    @("3425", "aaaasd", "123", "23",  "3123", "asda"),
    @("1231", "asdad",  "42",  "123", "24",   "asdasdas")
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $r = 2 + $i
    $rowValues = $data[$i]
    for ($j = 0; $j -lt $rowValues.Length; $j++) {
        $c = 1 + $j
        $ws.Cells.Item($r, $c).Value = $rowValues[$j]
    }
}
